$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("confusion")

# Updated confusion matrix values (leave-one-out update)
$values = @{
    8  = @(2181, 1896, 93, 69, 10, 6, 0)
    9  = @(1582, 4110, 267, 173, 32, 15, 2)
    10 = @(165, 652, 154, 82, 16, 8, 1)
    11 = @(64, 317, 83, 180, 19, 4, 0)
    12 = @(16, 57, 22, 35, 65, 4, 0)
    13 = @(5, 33, 12, 21, 5, 54, 0)
    14 = @(0, 1, 5, 2, 0, 1, 29)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        # columns E..K are columns 5..11
        $col = 5 + $i
        $ws.Cells.Item($row, $col).Value = $rowValues[$i]
    }
}

# Restore the active selection left behind by the editor
[void]$ws.Range("M19:N20").Select()
